# Hindalco prices sheet update (2025-10-10): insert the newest price record
# as row 2 (pushing the existing 41 rows of history down by one row), and
# rebuild the hyperlinks in column F so they point at the right rows again.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Insert a new blank row right under the header, shifting every
#    existing data row (old row 2..42) down to (new row 3..43).
# ---------------------------------------------------------------------
$ws.Range("A2:F2").EntireRow.Insert()

# ---------------------------------------------------------------------
# 2) Populate the new row 2 with the latest circular's data.
#    Column E holds a dd.mm.yyyy value that must stay plain text (like
#    every other row), so force a text format before typing it in.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = 42
$ws.Range("B2").Value = "2. P0610 (99.85% min) /P1020/ EC Grade Ingot & Sow 99.7% (min) / Cast Bar"
$ws.Range("C2").Value = "P1020"
$ws.Range("D2").Value = 282.25
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "10.10.2025"
$ws.Range("F2").Value = "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-10-october-2025.pdf"

# Re-apply the same look-and-feel the other data rows use (the row
# insert above borrowed the header's bold style) by pulling the format
# from row 3, now that every value is in place.
$ws.Range("A3:F3").Copy() | Out-Null
$ws.Range("A2:F2").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------
# 3) Rebuild column F's hyperlinks. The row insert moved the cell text
#    down with the rows, but left the old Hyperlinks collection pointing
#    at its original (now stale) rows, so drop them all and re-add them
#    in top-to-bottom order: the brand-new F2 link, then each existing
#    link shifted one row down (old F2->new F3, ... old F22->new F23).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Cells.Item(2, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-10-october-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(3, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-09-october-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(4, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-04-october-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(5, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-30-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(6, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(7, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(8, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-25-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(9, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(10, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(11, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-18-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(12, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-17-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(13, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(14, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(15, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-02-september-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(16, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-27-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(17, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-26-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(18, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-23-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(19, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-20-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(20, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-19-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(21, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-14-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(22, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-13-august-2025.pdf") | Out-Null
$ws.Hyperlinks.Add($ws.Cells.Item(23, 6), "https://www.hindalco.com/Upload/PDF/primary-ready-reckoner-12-august-2025.pdf") | Out-Null

# Adding a Hyperlinks entry auto-applies Excel's blue/underlined "Hyperlink"
# style, but this sheet has always kept column F in the plain data style
# (s=3, same as every other column) - so restore that look on F2:F23.
$ws.Range("A3:A3").Copy() | Out-Null
$ws.Range("F2:F23").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

"Updated Hindalco price sheet: inserted 10.10.2025 record, now " + $ws.UsedRange.Rows.Count + " rows."
